$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert three new date columns before column B. This shifts the
#    existing Jun_17 / Jun_15 / Jun_13 / Jun_10 columns (B:E) to the
#    right (E:H) while carrying their values/styles along automatically.
$ws.Range("B:D").EntireColumn.Insert()

# Carry the existing "8 characters wide" custom width forward onto the
# newly inserted date columns (C:D) and the columns that got shifted
# right by the insert (E:H); the brand-new "newest date" column (B)
# keeps the default width, same as before the edit.
$ws.Range("C:H").ColumnWidth = 7.14

# 2) Populate the header row with the two new report dates.
#    Jun_26 legitimately appears twice (columns C & D), matching the
#    source report, and is introduced before Jun_27 to mirror the
#    order in which the underlying report data was generated.
$ws.Range("C1").Value = "Jun_26"
$ws.Range("D1").Value = "Jun_26"

# 3) The newly inserted B:D columns are blank for every existing broker
#    row - fill them with the same "UN" (unchanged) placeholder used
#    throughout the rest of the sheet.
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 2).Value = "UN"
    $ws.Cells.Item($r, 3).Value = "UN"
    $ws.Cells.Item($r, 4).Value = "UN"
}

# 4) Wells Fargo & Co (row 21) got a new rating action on 6/19/2018 that
#    lands in both of the new Jun_26 columns.
$ws.Cells.Item(21, 3).Value = "6/19/2018,Reiterates,Buy,`$53.00"
$ws.Cells.Item(21, 4).Value = "6/19/2018,Reiterates,Buy,`$53.00"

# 5) Add the new broker group at the bottom of the sheet.
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28").Value = "UN"
$ws.Range("C28").Value = "UN"
$ws.Range("D28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29").Value = "UN"
$ws.Range("C29").Value = "UN"
$ws.Range("D29").Value = "UN"

# 6) Finally, the newest date column.
$ws.Range("B1").Value = "Jun_27"
